# Add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right before the existing "2022-Q2"
#    sheet and fill it with the quarterly fund-holdings table.
# 2. Insert a new summary row into "总计" (right after the header row) for
#    the 2022-Q3 totals, pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: new "2022-Q3" sheet, positioned before "2022-Q2"
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

# Header row (bold, matches the other quarter sheets)
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"
$q3.Range("B1:H1").Font.Bold = $true
$q3.Range("B1:H1").HorizontalAlignment = -4108
$q3.Range("B1:H1").VerticalAlignment = -4160
$q3.Range("B1:H1").Borders.LineStyle = 1

# Column A is a plain numeric row index (0-based), like on every other
# quarter sheet.
$q3.Range("A2").Value = 0
$q3.Range("A3").Value = 1
$q3.Range("A4").Value = 2
$q3.Range("A5").Value = 3
$q3.Range("A6").Value = 4
$q3.Range("A7").Value = 5

# Columns B-G hold text (fund codes keep leading zeros, the numeric-looking
# figures keep their original formatted text) - force text storage so
# values round-trip exactly instead of being parsed as numbers.
$q3.Range("B2:G7").NumberFormat = "@"

$q3.Range("B2").Value = "460005"
$q3.Range("C2").Value = "华泰柏瑞价值增长混合A"
$q3.Range("D2").Value = "9.00"
$q3.Range("E2").Value = "93.35"
$q3.Range("F2").Value = "1.97"
$q3.Range("G2").Value = "0.1773"
$q3.Range("H2").Value = 10

$q3.Range("B3").Value = "014158"
$q3.Range("C3").Value = "博时浦惠一年持有期混合A"
$q3.Range("D3").Value = "3.87"
$q3.Range("E3").Value = "48.81"
$q3.Range("F3").Value = "3.46"
$q3.Range("G3").Value = "0.1339"
$q3.Range("H3").Value = 5

$q3.Range("B4").Value = "014159"
$q3.Range("C4").Value = "博时浦惠一年持有期混合C"
$q3.Range("D4").Value = "0.36"
$q3.Range("E4").Value = "48.81"
$q3.Range("F4").Value = "3.46"
$q3.Range("G4").Value = "0.0125"
$q3.Range("H4").Value = 5

$q3.Range("B5").Value = "010663"
$q3.Range("C5").Value = "长江均衡成长混合A"
$q3.Range("D5").Value = "0.20"
$q3.Range("E5").Value = "81.80"
$q3.Range("F5").Value = "4.52"
$q3.Range("G5").Value = "0.0090"
$q3.Range("H5").Value = 2

$q3.Range("B6").Value = "010664"
$q3.Range("C6").Value = "长江均衡成长混合C"
$q3.Range("D6").Value = "0.05"
$q3.Range("E6").Value = "81.80"
$q3.Range("F6").Value = "4.52"
$q3.Range("G6").Value = "0.0023"
$q3.Range("H6").Value = 2

$q3.Range("B7").Value = "010037"
$q3.Range("C7").Value = "华泰柏瑞价值增长混合C"
$q3.Range("D7").Value = "0.05"
$q3.Range("E7").Value = "93.35"
$q3.Range("F7").Value = "1.97"
$q3.Range("G7").Value = "0.0010"
$q3.Range("H7").Value = 10

# Column A keeps the same centered/bold/bordered style used on column A of
# every other quarter sheet (copy the already-styled header cell's format
# onto the A2:A7 block, one source cell -> many destination cells).
$q3.Range("B1").Copy() | Out-Null
$q3.Range("A2:A7").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# Step 2: insert the new 2022-Q3 summary row into "总计"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Re-apply the formatting of the (now pushed-down) old row 2 so the new
# row matches the rest of the table (column A bold/centered/bordered,
# columns B-D unstyled).
$total.Range("A3:D3").Copy() | Out-Null
$total.Range("A2:D2").PasteSpecial(-4122) | Out-Null

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 0.34
